$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F15").Value = "93_referral_statement"
$ws.Range("F23").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F24").Value = "ppe"
$ws.Range("F25").Value = "ppe"
$ws.Range("F28").Value = "off target movement || application instructions || env warning - species || env warning - water"
$ws.Range("F33").Value = "135_product_information"
$ws.Range("F34").Value = "use restrictions"
$ws.Range("F37").Value = "mixing"
$ws.Range("F39").Value = "mixing"
$ws.Range("F40").Value = "mixing"
$ws.Range("F131").Value = "use restrictions"
$ws.Range("F133").Value = "154_pesticide_storage"
